# The deck ships two theme parts: the slide master's theme ("Integral" /
# "Red Violet" colour scheme) and the notes master's theme ("Office Theme").
# The authored edit swaps the two themes' content so the slide master now
# renders with the standard "Office" colour scheme instead of "Red Violet".
#
# The PowerPoint object model exposes theme colours through
# Master.Theme.ThemeColorScheme.Colors(index).RGB (index 1-12, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order) and COM RGB values
# are packed little-endian as 0xBBGGRR, so build each value from its R/G/B
# bytes instead of hand-transposing hex strings.

function ComRgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target colour scheme ("Office"), in slot order.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1      000000
    @(0xFF, 0xFF, 0xFF),  # 2  lt1      FFFFFF
    @(0x44, 0x54, 0x6A),  # 3  dk2      44546A
    @(0xE7, 0xE6, 0xE6),  # 4  lt2      E7E6E6
    @(0x5B, 0x9B, 0xD5),  # 5  accent1  5B9BD5
    @(0xED, 0x7D, 0x31),  # 6  accent2  ED7D31
    @(0xA5, 0xA5, 0xA5),  # 7  accent3  A5A5A5
    @(0xFF, 0xC0, 0x00),  # 8  accent4  FFC000
    @(0x44, 0x72, 0xC4),  # 9  accent5  4472C4
    @(0x70, 0xAD, 0x47),  # 10 accent6  70AD47
    @(0x05, 0x63, 0xC1),  # 11 hlink    0563C1
    @(0x95, 0x4F, 0x72)   # 12 folHlink 954F72
)

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $colorScheme.Colors($i).RGB = ComRgb $rgb[0] $rgb[1] $rgb[2]
}

# Best-effort: also rename the theme/colour-scheme to match the swapped-in
# "Office" theme (no-op on hosts that treat these as read-only).
$theme.Name = "Office Theme"
$colorScheme.Name = "Office"
